$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "w1p1"
$ws.Range("D3").Value = "w2p1"

$ws.Range("D4").Select()
